# Apply the cryptos price/volume update described in the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '29.634.08'
$ws.Range('E2').Value = '  +4.05%  '

$ws.Range('D3').Value = '1.916.47'

$ws.Range('D4').Value = "'1.002"
$ws.Range('E4').Value = '  +0.12%  '

$ws.Range('D5').Value = "'334.40"
$ws.Range('E5').Value = '  +0.98%  '

$ws.Range('E6').Value = '  +0.07%  '

$ws.Range('D7').Value = "'0.4678"
$ws.Range('E7').Value = '  +1.78%  '

$ws.Range('D8').Value = "'0.4114"
$ws.Range('E8').Value = '  +1.81%  '

$ws.Range('D9').Value = "'48.14"
$ws.Range('E9').Value = '  +1.18%  '

$ws.Range('D10').Value = "'0.08036"
$ws.Range('E10').Value = '  +2.20%  '

$ws.Range('D11').Value = "'1.013"
$ws.Range('E11').Value = '  +2.51%  '

$ws.Range('E12').Value = '  +4.36%  '

$ws.Range('D13').Value = '1.886.53'
$ws.Range('E13').Value = '  +1.26%  '

$ws.Range('D14').Value = "'5.988"
$ws.Range('E14').Value = '  +2.36%  '

$ws.Range('D15').Value = "'7.174"
$ws.Range('E15').Value = '  +2.20%  '

$ws.Range('D16').Value = "'89.96"
$ws.Range('E16').Value = '  +1.69%  '

$ws.Range('D17').Value = "'1.001"
$ws.Range('E17').Value = '  +0.18%  '

$ws.Range('E18').Value = '  +1.60%  '

$ws.Range('D19').Value = "'0.06585"
$ws.Range('E19').Value = '  +0.56%  '

$ws.Range('D20').Value = "'17.85"
$ws.Range('E20').Value = '  +3.69%  '

$ws.Range('E21').Value = '  +0.19%  '

$ws.Range('D22').Value = '29.618.17'
$ws.Range('E22').Value = '  +4.05%  '

$ws.Range('D23').Value = "'5.583"
$ws.Range('E23').Value = '  +4.54%  '

$ws.Range('D24').Value = "'11.59"
$ws.Range('E24').Value = '  +6.60%  '

$ws.Range('E25').Value = '  -1.72%  '

$ws.Range('D26').Value = '2.128.43'
$ws.Range('E26').Value = '  +2.06%  '

$ws.Range('D27').Value = "'155.51"
$ws.Range('E27').Value = '  -1.23%  '

$ws.Range('D28').Value = "'19.90"
$ws.Range('E28').Value = '  +3.16%  '

$ws.Range('D29').Value = "'5.753"
$ws.Range('E29').Value = '  +8.16%  '

$ws.Range('D30').Value = "'2.141"
$ws.Range('E30').Value = '  +3.38%  '

$ws.Range('D31').Value = "'117.57"
$ws.Range('E31').Value = '  +0.12%  '

$ws.Range('D32').Value = "'1.070"
$ws.Range('E32').Value = '  +11.66%  '

$ws.Range('D33').Value = "'0.09460"
$ws.Range('E33').Value = '  +1.31%  '

$ws.Range('D34').Value = "'1.430"
$ws.Range('E34').Value = '  +2.49%  '

$ws.Range('E35').Value = '  -0.66%  '

$ws.Range('D36').Value = "'5.403"
$ws.Range('E36').Value = '  +3.27%  '

$ws.Range('D37').Value = "'0.06127"
$ws.Range('E37').Value = '  +1.63%  '

$ws.Range('E38').Value = '  +2.76%  '

$ws.Range('D39').Value = "'8.414"
$ws.Range('E39').Value = '  +1.39%  '

$ws.Range('E40').Value = '  +1.07%  '

$ws.Range('D41').Value = "'0.5903"
$ws.Range('E41').Value = '  +2.22%  '

$ws.Range('D42').Value = "'0.1845"
$ws.Range('E42').Value = '  +1.80%  '

$ws.Range('D43').Value = "'10.21"
$ws.Range('E43').Value = '  +1.53%  '

$ws.Range('D44').Value = "'1.270"
$ws.Range('E44').Value = '  +2.34%  '

$ws.Range('D45').Value = "'2.345"
$ws.Range('E45').Value = '  +1.68%  '

$ws.Range('D46').Value = "'0.07519"
$ws.Range('E46').Value = '  +4.69%  '

$ws.Range('B47').Value = 'EnergySwap'
$ws.Range('C47').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D47').Value = "'12.23"
$ws.Range('E47').Value = '  +3.45%  '

$ws.Range('B48').Value = 'Decentraland'
$ws.Range('C48').Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range('D48').Value = "'0.5570"
$ws.Range('E48').Value = '  +2.43%  '

$ws.Range('D49').Value = "'1.932"
$ws.Range('E49').Value = '  +2.59%  '

$ws.Range('D50').Value = "'113.58"
$ws.Range('E50').Value = '  +3.33%  '

$ws.Range('D51').Value = "'0.2983"
$ws.Range('E51').Value = '  +11.68%  '
